$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text formatting (avoid Excel auto-converting
# numeric-looking strings into actual numbers and dropping trailing zeros, etc.)
$cells = @('D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E5', 'D6', 'E6', 'D8', 'E8', 'E9', 'D10', 'E10', 'E11', 'D12', 'E12', 'D14', 'E14', 'D15', 'E15', 'E16', 'D17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'E25', 'D26', 'E26', 'E27', 'D28', 'D29', 'E29', 'E30', 'E31', 'E32', 'E33', 'D34', 'E34', 'D35', 'E35', 'E36', 'E37', 'D38', 'E38', 'E39', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'E44', 'D45', 'E45', 'E46', 'D47', 'E47', 'E48', 'E49', 'D50', 'E50', 'E51')
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '67.228.66'
$ws.Range('E2').Value = '  +0.93%  '
$ws.Range('D3').Value = '3.512.32'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '596.57'
$ws.Range('E5').Value = '  +0.95%  '
$ws.Range('D6').Value = '173.42'
$ws.Range('E6').Value = '  +1.86%  '
$ws.Range('D8').Value = '0.592'
$ws.Range('E8').Value = '  +2.19%  '
$ws.Range('E9').Value = '  +5.84%  '
$ws.Range('D10').Value = '7.29'
$ws.Range('E10').Value = '  -0.72%  '
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '4.123.59'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D14').Value = '29.26'
$ws.Range('E14').Value = '  +3.48%  '
$ws.Range('D15').Value = '67.171.54'
$ws.Range('E15').Value = '  +0.79%  '
$ws.Range('E16').Value = '  +1.12%  '
$ws.Range('D17').Value = '3.510.21'
$ws.Range('D18').Value = '6.34'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('D19').Value = '14.16'
$ws.Range('E19').Value = '  +0.94%  '
$ws.Range('D20').Value = '395.99'
$ws.Range('E20').Value = '  +1.91%  '
$ws.Range('D21').Value = '8.02'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').Value = '73.15'
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').Value = '0.538'
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('E25').Value = '  +0.48%  '
$ws.Range('D26').Value = '10.27'
$ws.Range('E26').Value = '  -1.66%  '
$ws.Range('E27').Value = '  +1.23%  '
$ws.Range('D28').Value = '0.998'
$ws.Range('D29').Value = '6.33'
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('E30').Value = '  -0.95%  '
$ws.Range('E31').Value = '  +0.43%  '
$ws.Range('E32').Value = '  +1.28%  '
$ws.Range('E33').Value = '  -0.60%  '
$ws.Range('D34').Value = '1.68'
$ws.Range('E34').Value = '  +4.13%  '
$ws.Range('D35').Value = '163.27'
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('E36').Value = '  +0.87%  '
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('D38').Value = '7.06'
$ws.Range('E38').Value = '  +6.30%  '
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('D40').Value = '0.0752'
$ws.Range('E40').Value = '  +0.70%  '
$ws.Range('D41').Value = '26.60'
$ws.Range('E41').Value = '  +0.85%  '
$ws.Range('D42').Value = '27.22'
$ws.Range('E42').Value = '  +2.19%  '
$ws.Range('D43').Value = '2.837.14'
$ws.Range('E43').Value = '  +0.40%  '
$ws.Range('E44').Value = '  +3.31%  '
$ws.Range('D45').Value = '42.99'
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('E46').Value = '  -2.20%  '
$ws.Range('D47').Value = '338.40'
$ws.Range('E47').Value = '  -4.67%  '
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('E49').Value = '  +2.73%  '
$ws.Range('D50').Value = '6.50'
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('E51').Value = '  -0.43%  '
